$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# Append " [years]" to the label (column C) of every "Age at diagnosis of ..."
# / "Age at time of death" variable, matching the author's commit:
# "updated Dataschema P2 to include "[years]" in all AGE related variables
#  and DPEs accordingly"
$rows = 41, 43, 45, 47, 49, 51, 53, 55, 57, 60, 62

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 3)
    $label = [string]$cell.Value2
    if (-not $label.Contains("[years]")) {
        $cell.Value = $label + " [years]"
    }
}
